$d = $word.ActiveDocument
$wdStyleTypeCharacter = 2

# ---------------------------------------------------------------
# 1. Add the three new character styles used by the new paragraphs
# ---------------------------------------------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", $wdStyleTypeCharacter)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", $wdStyleTypeCharacter)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", $wdStyleTypeCharacter)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------
# Helper: apply a character style to every run whose text equals
# $searchText exactly (searching the whole document, repeatedly).
# ---------------------------------------------------------------
function Apply-StyleToAllMatches($styleName, $searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $rng.Style = $styleName
        $rng.Collapse(0)
    }
}

# ---------------------------------------------------------------
# Helper: locate a paragraph via a (possibly partial) anchor text,
# then apply a character style to the paragraph's *entire* text
# (excluding the trailing paragraph mark), so the whole run(s) of
# text get the style, not just the anchor substring.
# ---------------------------------------------------------------
function Apply-StyleToWholeParagraph($styleName, $anchorText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    if (-not $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        return
    }
    $para = $rng.Paragraphs(1).Range
    $para.MoveEnd(1, -1) | Out-Null
    $para.Style = $styleName
}

# ---------------------------------------------------------------
# 2. Apply GaNStyle to all four "Datas das campanhas..." runs
# ---------------------------------------------------------------
Apply-StyleToAllMatches "GaNStyle" "Datas das campanhas de 2022 que usam constelação de botas: 14 a 23 de maio, 13 a 22 de junho, 12 a 21 de julho"

# ---------------------------------------------------------------
# 3. Apply GaNParagraph to the "Está a participar..." paragraph run
# ---------------------------------------------------------------
Apply-StyleToWholeParagraph "GaNParagraph" "Está a participar numa campanha global"

# ---------------------------------------------------------------
# 4. Apply GaNLinks to the "por Jenik Hollan..." paragraph run
# ---------------------------------------------------------------
Apply-StyleToWholeParagraph "GaNLinks" "por Jenik Hollan, CzechGlobe"

Write-Host "Done applying GaN styles"
